$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D retains text formatting so numeric-looking price strings
# (e.g. "580.98") are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '69.342.54'
$ws.Range("E2").Value = '  +1.60%  '
$ws.Range("D3").Value = '3.399.45'
$ws.Range("E3").Value = '  +1.38%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '580.98'
$ws.Range("E5").Value = '  -0.49%  '
$ws.Range("D6").Value = '178.88'
$ws.Range("E6").Value = '  +1.06%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("E8").Value = '  +0.42%  '
$ws.Range("E9").Value = '  +8.23%  '
$ws.Range("E10").Value = '  +0.88%  '
$ws.Range("D11").Value = '48.32'
$ws.Range("E11").Value = '  +0.52%  '
$ws.Range("E12").Value = '  +3.45%  '
$ws.Range("D13").Value = '682.09'
$ws.Range("E13").Value = '  -0.50%  '
$ws.Range("D14").Value = '3.945.75'
$ws.Range("E14").Value = '  +1.32%  '
$ws.Range("E15").Value = '  +2.04%  '
$ws.Range("D16").Value = '69.434.11'
$ws.Range("E16").Value = '  +1.71%  '
$ws.Range("D17").Value = '3.398.79'
$ws.Range("E17").Value = '  +1.42%  '
$ws.Range("E18").Value = '  +0.70%  '
$ws.Range("D19").Value = '17.68'
$ws.Range("E19").Value = '  +1.34%  '
$ws.Range("D20").Value = '11.30'
$ws.Range("E20").Value = '  +0.73%  '
$ws.Range("D22").Value = '5.37'
$ws.Range("E22").Value = '  -1.61%  '
$ws.Range("D23").Value = '17.07'
$ws.Range("E23").Value = '  +0.56%  '
$ws.Range("D24").Value = '101.17'
$ws.Range("E24").Value = '  +0.70%  '
$ws.Range("E25").Value = '  -0.59%  '
$ws.Range("D26").Value = '2.71'
$ws.Range("E26").Value = '  +0.28%  '
$ws.Range("E27").Value = '  +2.24%  '
$ws.Range("D28").Value = '33.48'
$ws.Range("E28").Value = '  +1.45%  '
$ws.Range("D29").Value = '8.75'
$ws.Range("E29").Value = '  +2.65%  '
$ws.Range("E30").Value = '  -0.33%  '
$ws.Range("E31").Value = '  +13.51%  '
$ws.Range("D32").Value = '555.30'
$ws.Range("E32").Value = '  -0.25%  '
$ws.Range("D33").Value = '11.03'
$ws.Range("E33").Value = '  -0.50%  '
$ws.Range("E34").Value = '  +0.12%  '
$ws.Range("D35").Value = '58.07'
$ws.Range("E35").Value = '  +0.07%  '
$ws.Range("E36").Value = '  +0.06%  '
$ws.Range("D37").Value = '3.607.82'
$ws.Range("E37").Value = '  -2.94%  '
$ws.Range("E38").Value = '  +2.29%  '
$ws.Range("D39").Value = '35.28'
$ws.Range("E39").Value = '  +1.26%  '
$ws.Range("E40").Value = '  +10.98%  '
$ws.Range("E41").Value = '  +4.60%  '
$ws.Range("E42").Value = '  +3.48%  '
$ws.Range("D43").Value = '3.40'
$ws.Range("E43").Value = '  +5.07%  '
$ws.Range("E44").Value = '  +3.73%  '
$ws.Range("D45").Value = '0.336'
$ws.Range("E45").Value = '  +0.04%  '
$ws.Range("E46").Value = '  +1.17%  '
$ws.Range("E47").Value = '  +0.28%  '
$ws.Range("E49").Value = '  -0.19%  '
$ws.Range("D50").Value = '131.27'
$ws.Range("E50").Value = '  -0.50%  '
$ws.Range("D51").Value = '2.62'
$ws.Range("E51").Value = '  +1.72%  '

# Restore default (Normal) style on column D so no stray number-format
# style index is left on the cells, matching the original workbook look.
$ws.Range("D2:D51").Style = "Normal"

